# Update "想去人数" (want-to-go count) values that changed between scrapes.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 9880
$ws1.Range("F10").Value = 3943
$ws1.Range("F12").Value = 117
$ws1.Range("F16").Value = 557

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 9880
$ws4.Range("F11").Value = 3943
$ws4.Range("F13").Value = 117
$ws4.Range("F17").Value = 557
